$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.947.76"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.888.43"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7724"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3092"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07120"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08560"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.02%  "
$ws.Range("D12").Value = "1.993.02"
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7626"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.322"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.151"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "30.088.41"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.305.34"
$ws.Range("E19").Value = "  +8.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007762"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9986"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.958"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.296"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.028"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.432"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.540"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.495"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05432"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7461"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01955"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4457"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "1.106.58"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.077"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8474"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.198.55"
$ws.Range("E46").Value = "  +7.77%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.869"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.611"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
